# Update "想去人数" (want-to-go count) values in column F on sheets
# "展览" (sheet 1) and "全部类型" (sheet 4). Each sheet shares the same
# set of row values for these particular events, so the same deltas are
# applied in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value mapping for the "展览" sheet (column F)
$exhibitUpdates = @{
    3  = 551
    6  = 503
    7  = 102
    8  = 118
    10 = 6717
    12 = 370
    13 = 2997
    14 = 194
    15 = 337
    17 = 538
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value mapping for the "全部类型" sheet (column F)
$allTypeUpdates = @{
    5  = 551
    8  = 503
    9  = 102
    10 = 118
    13 = 6717
    16 = 370
    17 = 2997
    18 = 194
    19 = 337
    21 = 538
}

foreach ($row in $allTypeUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypeUpdates[$row]
}
